# Add new keyword rows to the Cluster_Keywords table, re-sort the table
# (Category then Stem, matching its existing saved sort order), and add a
# conditional-formatting rule that highlights duplicate Stem values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item("Cluster_Keywords")

$newRows = @(
    @("Crech", "Childcare"),
    @("Kinde", "Childcare"),
    @("Chemi", "Healthcare"),
    @("DHL", "Logistics"),
    @("Truck", "Logistics"),
    @("Whole", "Logistics")
)

foreach ($pair in $newRows) {
    $stem = $pair[0]
    $category = $pair[1]

    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

    $row = $tbl.ListRows.Add()
    $r = $row.Range.Row

    $ws.Cells.Item($r, 1).Value = $stem
    $ws.Cells.Item($r, 2).Formula = "=LEN(Cluster_Keywords[[#This Row],[Stem]])"
    $ws.Cells.Item($r, 3).Value = $category

    $fmtSrc = $ws.Range($ws.Cells.Item($lastRow, 1), $ws.Cells.Item($lastRow, 3))
    $fmtDst = $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 3))
    $fmtSrc.Copy()
    $fmtDst.PasteSpecial(-4122)
}

$dataRange = $tbl.DataBodyRange
$lastDataRow = $dataRange.Row + $dataRange.Rows.Count - 1

[void]$tbl.Sort.SortFields.Clear()
[void]$tbl.Sort.SortFields.Add($ws.Range("C2:C" + $lastDataRow))
[void]$tbl.Sort.SortFields.Add($ws.Range("A2:A" + $lastDataRow))
$tbl.Sort.Header = 1
[void]$tbl.Sort.Apply()

$dupRange = $ws.Range("A2:A" + $lastDataRow)
$cf = $dupRange.FormatConditions.AddUniqueValues()
$cf.DupeUnique = 1
$cf.Font.Color = 393372
$cf.Interior.Color = 13551615

$ws.Cells.Item($lastDataRow + 1, 1).Select()

"Added " + $newRows.Count + " rows; table now spans to row " + $lastDataRow
